$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A48").Value = "Golang Developer"
$ws.Range("B48").Value = "https://www.dice.com/job-detail/660f87c2-8f1e-4464-96fc-764482a91acd"
$ws.Range("C48").Value = "Remote"
$ws.Range("D48").Value = "Contract"
$ws.Range("E48").Value = "Depends on Experience"
$ws.Range("F48").Value = "Montek System"
